$d = $word.ActiveDocument

# 1. Remove the empty paragraph that only contains the horizontal-rule
#    VML <w:pict> placed right after the "Thomas J. Herzog" heading
#    (the second horizontal rule, before "Education and Certificates",
#    stays untouched).
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $para = $d.Paragraphs.Item($i)
    $ptext = $para.Range.Text
    # an empty paragraph (just the pict run) has no visible text other
    # than the paragraph mark
    if ($ptext -eq "`r" -or $ptext -eq "") {
        $prev = $null
        if ($i -gt 1) { $prev = $d.Paragraphs.Item($i - 1) }
        if ($prev -ne $null -and $prev.Range.Text -match "Thomas J\. Herzog") {
            $para.Range.Delete()
            break
        }
    }
}

# 2. Skills line: swap the tech list.
$d.Content.Find.Execute(
    "Knex • PostgreSQL • Git • Github • Netlify • Docker • NPM • Yarn • Babel • Webpack • Eslint •",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "MySQL • PostgreSQL • Git • Github • Docker • NPM • Yarn • Babel • Webpack • Eslint • PHP •",
    2) | Out-Null

# 3. Experience section updates.

# 2019 - Present -> 2020 - Present
$d.Content.Find.Execute(
    "2019 - Present", $true, $false, $false, $false, $false, $true, 1, $false,
    "2020 - Present", 2) | Out-Null

# Job title / employer line
$d.Content.Find.Execute(
    "Produce Lead - New Seasons Market", $true, $false, $false, $false, $false, $true, 1, $false,
    "Fullstack software Developer - Freelance", 2) | Out-Null

# Description of the (now freelance) role
$d.Content.Find.Execute(
    "Training staff and overseeing their productivity, ensuring quality standards are met.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Working with a team of developers to enhance and maintain a custom-built internal ERP business tool. • JavaScript • HTML • PHP • CSS • Docker • MySQL",
    2) | Out-Null

# 2018 - 2019 -> 2018 - Present
$d.Content.Find.Execute(
    "2018 - 2019", $true, $false, $false, $false, $false, $true, 1, $false,
    "2018 - Present", 2) | Out-Null

# Produce Clerk title update
$d.Content.Find.Execute(
    "Produce Clerk - New Seasons Market", $true, $false, $false, $false, $false, $true, 1, $false,
    "Produce Clerk / Lead - New Seasons Market", 2) | Out-Null

# Produce Clerk description (now reuses the old Produce Lead description)
$d.Content.Find.Execute(
    "Provided service based on customer needs and complete tasks in a timely manner.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Training staff and overseeing their productivity, ensuring quality standards are met.",
    2) | Out-Null
